$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C29 has been removed from the design. Row 5 ("100nF" capacitors) held C29
# among its references - update the reference list, quantity, and add a
# comment noting its removal. Set the comment (I5) first and the reference
# list (B5) second so the new shared strings land in the same order as the
# authoritative edit (comment string, then the trimmed reference list).
$ws.Range("I5").Value = "C29 no longer exists"
$ws.Range("B5").Value = "C1 C12 C14 C15 C16 C18 C21 C23 C24 C26 C27 C3 C30 C32 C33 C36 C37 C38 C39 C5 C6 C8 C9 "
$ws.Range("C5").Value = 23

# Move the active selection to I8, matching where the user's cursor ended up.
$ws.Range("I8").Select() | Out-Null
